# Updates cryptos list prices (col D) and Volume(1h) percentages (col E)
# with fresh scraped values; rows 33/34 (Filecoin/VeChain) swap order.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "52.091.11"
$ws.Range("E2").Value = "  -0.15%  "
$ws.Range("D3").Value = "2.841.85"
$ws.Range("E3").Value = "  +1.68%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "362.12"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +5.46%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "113.75"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -2.73%  "
$ws.Range("E7").Value = "  +4.98%  "
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.603"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +3.88%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "41.71"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -1.31%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0862"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -0.86%  "
$ws.Range("E12").Value = "  +1.28%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "20.03"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -0.72%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.78"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +1.52%  "
$ws.Range("D15").Value = "3.285.43"
$ws.Range("E15").Value = "  +1.70%  "
$ws.Range("D16").Value = "2.844.92"
$ws.Range("E16").Value = "  +0.68%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.914"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +2.85%  "
$ws.Range("D18").Value = "51.849.58"
$ws.Range("E18").Value = "  -0.33%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.47"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +6.93%  "
$ws.Range("E20").Value = "  -1.62%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.56"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +1.28%  "
$ws.Range("D22").Value = "0.0₃0993"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "70.26"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -0.06%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "267.15"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -3.97%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.84"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +0.13%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "27.20"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +1.12%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.00"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +0.03%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.45"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +2.52%  "
$ws.Range("E29").Value = "  +1.50%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "53.13"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +5.37%  "
$ws.Range("E31").Value = "  -1.91%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "34.14"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -1.75%  "
$ws.Range("B33").Value = "VeChain"
$ws.Range("C33").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0451"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +20.08%  "
$ws.Range("B34").Value = "Filecoin"
$ws.Range("C34").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.89"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +3.50%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.33"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +6.72%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0843"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +2.34%  "
$ws.Range("E37").Value = "  +0.12%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.31"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +0.74%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.09"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -2.44%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "18.35"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -3.43%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "24.08"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +2.25%  "
$ws.Range("E42").Value = "  +2.03%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "127.65"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -0.34%  "
$ws.Range("E44").Value = "  -7.27%  "
$ws.Range("E45").Value = "  -3.28%  "
$ws.Range("D46").Value = "2.123.54"
$ws.Range("E46").Value = "  +0.66%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.39"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +1.38%  "
$ws.Range("E48").Value = "  +1.05%  "
$ws.Range("E49").Value = "  +8.80%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "5.85"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +4.99%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "9.04"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +0.96%  "
